$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 1 title bar: extend style across E1:H1 and merge A1:H1
# ---------------------------------------------------------------------------
$ws.Range("A1:D1").Copy()
$ws.Range("E1:H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1:H1").Merge()

# ---------------------------------------------------------------------------
# 2. Row 15: convert the blank row into the "Sprint 2" title bar
# ---------------------------------------------------------------------------
$ws.Range("A1:D1").Copy()
$ws.Range("A15:H15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A15").Value = "Sprint 2"
$ws.Range("A15:H15").Merge()
$ws.Rows.Item(15).RowHeight = 23.25

# ---------------------------------------------------------------------------
# 3. Row 16: Sprint 2 header row (copy format from row 2)
# ---------------------------------------------------------------------------
$ws.Range("A2:H2").Copy()
$ws.Range("A16:H16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A16").Value = $ws.Range("A2").Value
$ws.Range("B16").Value = $ws.Range("B2").Value
$ws.Range("C16").Value = $ws.Range("A2").Value
$ws.Range("D16").Value = $ws.Range("B2").Value
$ws.Range("E16").Value = $ws.Range("E2").Value
$ws.Range("F16").Value = $ws.Range("F2").Value
$ws.Range("G16").Value = $ws.Range("G2").Value
$ws.Range("H16").Value = $ws.Range("H2").Value
$ws.Rows.Item(16).RowHeight = 20.25

# ---------------------------------------------------------------------------
# 4. Row 17 (style like row 3)
# ---------------------------------------------------------------------------
$ws.Range("A3:H3").Copy()
$ws.Range("A17:H17").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(17).RowHeight = 19.5
$ws.Range("A17").Value = 43789
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 43789
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 34
$ws.Range("F17").Formula = "=B17"
$ws.Range("G17").Formula = "=`$E`$3-F17"
$ws.Range("H17").Formula = "=D18"

# ---------------------------------------------------------------------------
# 5. Row 18 (style like row 4)
# ---------------------------------------------------------------------------
$ws.Range("A4:H4").Copy()
$ws.Range("A18:H18").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(18).RowHeight = 18.75
$ws.Range("A18").Value = 43790
$ws.Range("C18").Value = 43772
$ws.Range("D18").Formula = "=E17"
$ws.Range("F18").Formula = "=B18+F17"
$ws.Range("H18").Formula = "=D17"

# ---------------------------------------------------------------------------
# 6. Rows 19-24 (style like row 5)
# ---------------------------------------------------------------------------
$ws.Range("A5:H5").Copy()
$ws.Range("A19:H24").PasteSpecial(-4122)  # xlPasteFormats
for ($r = 19; $r -le 24; $r++) {
    $ws.Rows.Item($r).RowHeight = 18.75
    $ws.Range("A$r").Value = 43789 + ($r - 17)
    $ws.Range("B$r").ClearContents()
}

# ---------------------------------------------------------------------------
# 7. Rows 25-28 (style like row 11 : no H column)
# ---------------------------------------------------------------------------
$ws.Range("A11:G11").Copy()
$ws.Range("A25:G28").PasteSpecial(-4122)  # xlPasteFormats
for ($r = 25; $r -le 28; $r++) {
    $ws.Rows.Item($r).RowHeight = 18.75
    $ws.Range("A$r").Value = 43789 + ($r - 17)
    $ws.Range("H$r").Clear()
}

# ---------------------------------------------------------------------------
# 8. Rows 29-30 (style like row 11, but also drop C/D)
# ---------------------------------------------------------------------------
$ws.Range("A11:G11").Copy()
$ws.Range("A29:G30").PasteSpecial(-4122)  # xlPasteFormats
for ($r = 29; $r -le 30; $r++) {
    $ws.Rows.Item($r).RowHeight = 18.75
    $ws.Range("A$r").Value = 43789 + ($r - 17)
    $ws.Range("C$r").Clear()
    $ws.Range("D$r").Clear()
    $ws.Range("H$r").Clear()
}

# ---------------------------------------------------------------------------
# 9. Shared formulas for F19:F30 and G18:G30
# ---------------------------------------------------------------------------
$ws.Range("G18:G30").Formula = "=`$E`$3-F18"
$ws.Range("F19:F30").Formula = "=B19+F18"

# ---------------------------------------------------------------------------
# 10. Shrink the Sprint-1 shared-formula ranges to their real extent
# ---------------------------------------------------------------------------
$ws.Range("G4:G10").Formula = "=`$E`$3-F4"
$ws.Range("F5:F10").Formula = "=B5+F4"

# ---------------------------------------------------------------------------
# 11. Move / resize the burn-down chart
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$co.Top = 0.0000787402
$co.Left = 1037.91796875
$co.Width = 509.4911321973425
$co.Height = 268.595905511811

# ---------------------------------------------------------------------------
# 12. Selection
# ---------------------------------------------------------------------------
[void]$ws.Range("D17").Select()
